# The "locked" security property's value (row 3) on the "properties" sheet
# changes from "true" to "false".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

# Leading apostrophe forces Excel to store the value as text (shared string)
# instead of re-interpreting "false" as a boolean literal.
$ws.Range("E3").Value = "'false"

# Leave the active selection where it ended up after the edit.
$ws.Activate()
$ws.Range("E5").Select()
